$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 385 (shifts existing rows 385:414 down to 390:419)
$ws.Range("A385:T389").Insert(-4121)

# Common/constant column values for this data block (Cereza - Vega Modelo de Temuco)
$mercadoId = 10
$mercado = "Vega Modelo de Temuco"
$region = "La Araucanía"
$codreg = 9
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103001
$categoria = "Cereza"

# New row data: Date, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Origen, PrecioKg, KgUnidad
$newRows = @(
    @{ Row=385; D=45267; K="Bing";       L="Primera"; M=630; N=800;  O=1000; P=911;  Q="$/kilo (en caja de 15 kilos)"; R="Región del Maule";    S=911;  T=1 },
    @{ Row=386; D=45267; K="Bing";       L="Segunda"; M=250; N=600;  O=600;  P=600;  Q="$/kilo (en caja de 15 kilos)"; R="Región del Maule";    S=600;  T=1 },
    @{ Row=387; D=45267; K="Bing";       L="Tercera"; M=180; N=500;  O=500;  P=500;  Q="$/kilo (en caja de 15 kilos)"; R="Región del Maule";    S=500;  T=1 },
    @{ Row=388; D=45267; K="Brooks";     L="Primera"; M=380; N=800;  O=800;  P=800;  Q="$/kilo (en caja de 15 kilos)"; R="Región del Maule";    S=800;  T=1 },
    @{ Row=389; D=45267; K="Royal Dawn"; L="Primera"; M=560; N=1000; O=1200; P=1089; Q="$/kilo (en caja de 15 kilos)"; R="Región de O'Higgins"; S=1089; T=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
